$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7959397435188293
$ws.Range("B1").Value = 0.6540499925613403
$ws.Range("C1").Value = 0.5846021771430969
$ws.Range("D1").Value = 0.6453571915626526
$ws.Range("E1").Value = 0.8075680136680603
